# Auto-generated edit script applying the Sagittarius_Profits.xlsx diff
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets of the combined workbook)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 19
$ws.Range("I6").Value = 19
$ws.Range("K6").Value = 57
$ws.Range("M6").Value = 55
$ws.Range("H18").Value = 938
$ws.Range("I18").Value = 883.3333
$ws.Range("J18").Value = 1102
$ws.Range("K18").Value = 883.3333
$ws.Range("L18").Value = 1102
$ws.Range("M18").Value = -599.3333
$ws.Range("N18").Value = -1670
$ws.Range("H31").Value = 83333390
$ws.Range("I31").Value = 83333390
$ws.Range("K31").Value = 250000170
$ws.Range("M31").Value = -249999940
$ws.Range("H32").Value = 3810.4666
$ws.Range("I32").Value = 3120
$ws.Range("J32").Value = 4155.7
$ws.Range("K32").Value = 3120
$ws.Range("L32").Value = 4155.7
$ws.Range("M32").Value = -2794
$ws.Range("N32").Value = -4807.7
$ws.Range("H33").Value = 241.44444
$ws.Range("I33").Value = 166.2
$ws.Range("J33").Value = 335.5
$ws.Range("K33").Value = 166.2
$ws.Range("L33").Value = 335.5
$ws.Range("M33").Value = 62.80000000000001
$ws.Range("N33").Value = -793.5
$ws.Range("H38").Value = 973.7917
$ws.Range("I38").Value = 521.94116
$ws.Range("J38").Value = 2071.1428
$ws.Range("K38").Value = 1565.82348
$ws.Range("L38").Value = 6213.428400000001
$ws.Range("M38").Value = -1193.82348
$ws.Range("N38").Value = -6957.428400000001
$ws.Range("H64").Value = 5331.6665
$ws.Range("H67").Value = 5331.6665
$ws.Range("H70").Value = 66143.31
$ws.Range("J70").Value = 70399.53
$ws.Range("L70").Value = 211198.59
$ws.Range("N70").Value = -211738.59
$ws.Range("H73").Value = 66143.31
$ws.Range("J73").Value = 70399.53
$ws.Range("L73").Value = 211198.59
$ws.Range("N73").Value = -213070.59
$ws.Range("H137").Value = 1396.8422
$ws.Range("I137").Value = 1428
$ws.Range("J137").Value = 1309.6
$ws.Range("K137").Value = 4284
$ws.Range("L137").Value = 3928.8
$ws.Range("M137").Value = -1734
$ws.Range("N137").Value = -9028.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 2728
$ws.Range("I21").Value = 2728
$ws.Range("K21").Value = 2728
$ws.Range("M21").Value = -2354
$ws.Range("H45").Value = 4710.5454
$ws.Range("I45").Value = 4979.5557
$ws.Range("K45").Value = 4979.5557
$ws.Range("M45").Value = -4602.5557
$ws.Range("H63").Value = 35935
$ws.Range("I63").Value = 35935
$ws.Range("K63").Value = 35935
$ws.Range("M63").Value = -35249
$ws.Range("H66").Value = 35935
$ws.Range("I66").Value = 35935
$ws.Range("K66").Value = 179675
$ws.Range("M66").Value = -176243
$ws.Range("H74").Value = 1515.5
$ws.Range("I74").Value = 1323.25
$ws.Range("K74").Value = 1323.25
$ws.Range("M74").Value = -449.25
$ws.Range("H77").Value = 1515.5
$ws.Range("I77").Value = 1323.25
$ws.Range("K77").Value = 6616.25
$ws.Range("M77").Value = -2248.25
$ws.Range("H110").Value = 7400766
$ws.Range("I110").Value = 7400766
$ws.Range("K110").Value = 7400766
$ws.Range("M110").Value = -7398721
$ws.Range("H132").Value = 1883.4286
$ws.Range("I132").Value = 1883.4286
$ws.Range("K132").Value = 5650.2858
$ws.Range("M132").Value = -3120.2858
$ws.Range("H134").Value = 59142.668
$ws.Range("J134").Value = 59142.668
$ws.Range("L134").Value = 59142.668
$ws.Range("N134").Value = -69282.66800000001
$ws.Range("H135").Value = 510000
$ws.Range("J135").Value = 510000
$ws.Range("L135").Value = 510000
$ws.Range("N135").Value = -520140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1381.3334
$ws.Range("I20").Value = 1200
$ws.Range("K20").Value = 1200
$ws.Range("M20").Value = -953
$ws.Range("H33").Value = 23332
$ws.Range("J33").Value = 27500
$ws.Range("L33").Value = 27500
$ws.Range("N33").Value = -28172
$ws.Range("H86").Value = 1688.9
$ws.Range("J86").Value = 1526.6666
$ws.Range("L86").Value = 1526.6666
$ws.Range("N86").Value = -3772.6666
$ws.Range("H89").Value = 1688.9
$ws.Range("J89").Value = 1526.6666
$ws.Range("L89").Value = 7633.333000000001
$ws.Range("N89").Value = -18865.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2260.3333
$ws.Range("I31").Value = 1707
$ws.Range("K31").Value = 1707
$ws.Range("M31").Value = -1412
$ws.Range("H34").Value = 2260.3333
$ws.Range("I34").Value = 1707
$ws.Range("K34").Value = 1707
$ws.Range("M34").Value = -1505
$ws.Range("H58").Value = 2769.3333
$ws.Range("I58").Value = 2766.5454
$ws.Range("J58").Value = 2800
$ws.Range("K58").Value = 2766.5454
$ws.Range("L58").Value = 2800
$ws.Range("M58").Value = -2563.5454
$ws.Range("N58").Value = -3206
$ws.Range("H75").Value = 99999
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 99999
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 99999
$ws.Range("N75").Value = -101995
$ws.Range("H78").Value = 99999
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 99999
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 299997
$ws.Range("N78").Value = -309981
$ws.Range("H132").Value = 2312.6155
$ws.Range("I132").Value = 2246.4
$ws.Range("K132").Value = 6739.200000000001
$ws.Range("M132").Value = -4209.200000000001
$ws.Range("H136").Value = 2769.3333
$ws.Range("I136").Value = 2766.5454
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 8299.636200000001
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -5749.636200000001
$ws.Range("N136").Value = -13500
$ws.Range("M75").ClearContents()
$ws.Range("M78").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10001741
$ws.Range("I4").Value = 22000020
$ws.Range("J4").Value = 3174.3333
$ws.Range("K4").Value = 66000060
$ws.Range("L4").Value = 9522.999899999999
$ws.Range("M4").Value = -65999948
$ws.Range("N4").Value = -9746.999899999999
$ws.Range("H7").Value = 55
$ws.Range("I7").Value = 50
$ws.Range("K7").Value = 150
$ws.Range("M7").Value = -38
$ws.Range("H34").Value = 647.375
$ws.Range("J34").Value = 857.6
$ws.Range("L34").Value = 2572.8
$ws.Range("N34").Value = -2740.8
$ws.Range("H75").Value = 5117.8
$ws.Range("J75").Value = 7535
$ws.Range("L75").Value = 22605
$ws.Range("N75").Value = -24601
$ws.Range("H78").Value = 5117.8
$ws.Range("J78").Value = 7535
$ws.Range("L78").Value = 67815
$ws.Range("N78").Value = -77799
$ws.Range("H132").Value = 18999.666
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 18999.666
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 170996.994
$ws.Range("N132").Value = -176056.994
$ws.Range("H139").Value = 1584.5555
$ws.Range("I139").Value = 1504.875
$ws.Range("K139").Value = 4514.625
$ws.Range("M139").Value = 625.375
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 286.33334
$ws.Range("I97").Value = 265.8
$ws.Range("K97").Value = 265.8
$ws.Range("M97").Value = 230.2
$ws.Range("H132").Value = 949
$ws.Range("I132").Value = 949
$ws.Range("K132").Value = 2847
$ws.Range("M132").Value = -317

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2249.25
$ws.Range("I16").Value = 1998.3334
$ws.Range("J16").Value = 3002
$ws.Range("K16").Value = 1998.3334
$ws.Range("L16").Value = 3002
$ws.Range("M16").Value = -1828.3334
$ws.Range("N16").Value = -3342
$ws.Range("H40").Value = 2842.1428
$ws.Range("I40").Value = 2277.5386
$ws.Range("K40").Value = 2277.5386
$ws.Range("M40").Value = -2141.5386
$ws.Range("H82").Value = 1905.7368
$ws.Range("I82").Value = 2714.6365
$ws.Range("J82").Value = 793.5
$ws.Range("K82").Value = 2714.6365
$ws.Range("L82").Value = 793.5
$ws.Range("M82").Value = -2353.6365
$ws.Range("N82").Value = -1515.5
$ws.Range("H85").Value = 1905.7368
$ws.Range("I85").Value = 2714.6365
$ws.Range("J85").Value = 793.5
$ws.Range("K85").Value = 2714.6365
$ws.Range("L85").Value = 793.5
$ws.Range("M85").Value = -1466.6365
$ws.Range("N85").Value = -3289.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15235.9
$ws.Range("I62").Value = 16134
$ws.Range("J62").Value = 14851
$ws.Range("K62").Value = 16134
$ws.Range("L62").Value = 14851
$ws.Range("M62").Value = -15510
$ws.Range("N62").Value = -16099
$ws.Range("H65").Value = 15235.9
$ws.Range("I65").Value = 16134
$ws.Range("J65").Value = 14851
$ws.Range("K65").Value = 80670
$ws.Range("L65").Value = 74255
$ws.Range("M65").Value = -77550
$ws.Range("N65").Value = -80495
$ws.Range("H81").Value = 1114149
$ws.Range("I81").Value = 3973.1667
$ws.Range("K81").Value = 7946.3334
$ws.Range("M81").Value = -6885.3334
$ws.Range("H84").Value = 1114149
$ws.Range("I84").Value = 3973.1667
$ws.Range("K84").Value = 39731.667
$ws.Range("M84").Value = -34427.667
$ws.Range("H107").Value = 1347.5
$ws.Range("I107").Value = 1397.1666
$ws.Range("J107").Value = 1198.5
$ws.Range("K107").Value = 4191.4998
$ws.Range("L107").Value = 3595.5
$ws.Range("M107").Value = -2271.4998
$ws.Range("N107").Value = -7435.5
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 6000
$ws.Range("M113").Value = -3830
$ws.Range("H136").Value = 2797.3225
$ws.Range("I136").Value = 2720.92
$ws.Range("J136").Value = 3115.6667
$ws.Range("K136").Value = 8162.76
$ws.Range("L136").Value = 9347.000100000001
$ws.Range("M136").Value = -5612.76
$ws.Range("N136").Value = -14447.0001
